{"js": "const replacements = [\n  [\"915\u00f73=\", \"237\u00f79=\"],\n  [\"891\u00f72=\", \"884\u00f76=\"],\n  [\"740\u00f72=\", \"588\u00f76=\"],\n  [\"829\u00f79=\", \"733\u00f73=\"],\n  [\"293\u00f73=\", \"684\u00f78=\"],\n  [\"332\u00f78=\", \"591\u00f79=\"],\n  [\"546\u00f72=\", \"475\u00f79=\"],\n  [\"290\u00f77=\", \"922\u00f77=\"],\n  [\"857\u00f72=\", \"842\u00f73=\"],\n  [\"854\u00f73=\", \"826\u00f74=\"],\n  [\"376\u00f78=\", \"572\u00f76=\"],\n  [\"856\u00f78=\", \"526\u00f74=\"],\n  [\"576\u00f79=\", \"667\u00f79=\"],\n  [\"986\u00f79=\", \"749\u00f76=\"],\n  [\"660\u00f78=\", \"249\u00f77=\"],\n  [\"296\u00f76=\", \"503\u00f78=\"],\n  [\"868\u00f72=\", \"128\u00f73=\"],\n  [\"262\u00f72=\", \"123\u00f79=\"],\n  [\"850\u00f73=\", \"684\u00f77=\"],\n  [\"822\u00f77=\", \"443\u00f73=\"],\n  [\"966\u00f79=\", \"131\u00f78=\"],\n  [\"929\u00f73=\", \"949\u00f76=\"],\n  [\"807\u00f77=\", \"261\u00f73=\"],\n  [\"430\u00f78=\", \"929\u00f76=\"],\n  [\"994\u00f75=\", \"777\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"915\u00f73=\", \"237\u00f79=\"),\n    @(\"891\u00f72=\", \"884\u00f76=\"),\n    @(\"740\u00f72=\", \"588\u00f76=\"),\n    @(\"829\u00f79=\", \"733\u00f73=\"),\n    @(\"293\u00f73=\", \"684\u00f78=\"),\n    @(\"332\u00f78=\", \"591\u00f79=\"),\n    @(\"546\u00f72=\", \"475\u00f79=\"),\n    @(\"290\u00f77=\", \"922\u00f77=\"),\n    @(\"857\u00f72=\", \"842\u00f73=\"),\n    @(\"854\u00f73=\", \"826\u00f74=\"),\n    @(\"376\u00f78=\", \"572\u00f76=\"),\n    @(\"856\u00f78=\", \"526\u00f74=\"),\n    @(\"576\u00f79=\", \"667\u00f79=\"),\n    @(\"986\u00f79=\", \"749\u00f76=\"),\n    @(\"660\u00f78=\", \"249\u00f77=\"),\n    @(\"296\u00f76=\", \"503\u00f78=\"),\n    @(\"868\u00f72=\", \"128\u00f73=\"),\n    @(\"262\u00f72=\", \"123\u00f79=\"),\n    @(\"850\u00f73=\", \"684\u00f77=\"),\n    @(\"822\u00f77=\", \"443\u00f73=\"),\n    @(\"966\u00f79=\", \"131\u00f78=\"),\n    @(\"929\u00f73=\", \"949\u00f76=\"),\n    @(\"807\u00f77=\", \"261\u00f73=\"),\n    @(\"430\u00f78=\", \"929\u00f76=\"),\n    @(\"994\u00f75=\", \"777\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $old, $true, $false, $false, $false, $false, $true, 1, $false,\n        $new, 2\n    )\n}\n"}
